# valid login page for testcase1
#
# Rewrites the login-credentials worksheet: header labels and sample
# credentials are changed, and the sheet itself is renamed from "Test1"
# to "ValidLogin". The sheet view is also updated (zoom + selection).

$wb  = $excel.ActiveWorkbook
$old = $wb.Worksheets("Test1")

# Duplicate "Test1" in place. This gives the new sheet the next free
# sheetId (2, since "Test1" already holds 1) and - unlike
# Worksheets.Add() - carries over the worksheet's existing formatting /
# markup instead of starting from a blank sheet.
$old.Copy($old)
$new = $wb.Worksheets("Test1 (2)")

$new.Range("A1").Value = "UserName"
$new.Range("B1").Value = "Password"
$new.Range("A2").Value = "admin"
$new.Range("B2").Value = "manager"

$new.Name = "ValidLogin"

# Drop the original sheet now that its data has been recreated on the
# copy (re-resolve by name rather than reuse a cached reference, so we
# delete "Test1" and not the sheet we just renamed).
$wb.Worksheets("Test1").Delete()

# Match the saved view state: zoom 160% and selection on B1.
$ws = $wb.Worksheets("ValidLogin")
$ws.Select()
$excel.ActiveWindow.Zoom = 160
$ws.Range("B1").Select()
